$p = $ppt.ActivePresentation

# --- 1. Nudge the "Time step 2009" caption textbox on slide 1 down slightly ---
$s1 = $p.Slides.Item(1)
$caption = $s1.Shapes.Item(3)
$caption.Top = 338.55002

# --- 2. Add a new blank slide (3D route illustration placeholder) at the end ---
$newSlide = $p.Slides.Add($p.Slides.Count + 1, 7)
